# feat: (upload-service) remove hardcoded sheet names
#
# Rename the sheet currently named "English" to "language_English" and
# make it the active sheet (select/activate it), matching the behaviour
# a user gets in Excel when they rename and then click on the sheet tab.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("English")
$ws.Name = "language_English"
$ws.Activate()
